$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (headers) ---
$ws.Range("A1").Value = "privacyType"
$ws.Range("B1").Value = "post text"
$ws.Range("C1").Value = "expected result"
$ws.Range("D1").Value = "Result"

# --- Row 2 ---
$ws.Range("A2").Value = "Only me"
$ws.Range("B2").Value = "Helo only me post"
$ws.Range("C2").Value = "Shared with Only me"
$ws.Range("D2").Value = "Pass"

# --- Row 3 ---
$ws.Range("A3").Value = "Public"
$ws.Range("B3").Value = "Hello  Public post"
$ws.Range("C3").Value = "Shared with Public"
$ws.Range("D3").Value = "Pass"

# --- View state: zoom + selection on D2 ---
$excel.ActiveWindow.Zoom = 141
[void]$ws.Range("D2").Select()
